
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Resources" (sheet1) - add 3 new rows of resources (48-50)
# ---------------------------------------------------------------
$wsRes = $wb.Worksheets.Item("Resources")

# Row 48 - Interpretable Machine Learning book
$wsRes.Range("A48").Value = "Web"
$wsRes.Range("B48").Value = "Model"
$wsRes.Range("C48").Value = "Reference"
$wsRes.Range("D48").Value = "Mathematics"
$wsRes.Range("F48").Value = "Interpretable Machine Learning  A Guide for Making Black Box Models Explainable."
$wsRes.Range("F48").WrapText = $true
$wsRes.Range("G48").Value = "Christoph Molnar"
$wsRes.Range("H48").Value = "https://christophm.github.io/interpretable-ml-book/"
$wsRes.Rows.Item(48).RowHeight = 51

# Row 49 - Probability for Machine Learning
$wsRes.Range("A49").Value = "Web, Book"
$wsRes.Range("B49").Value = "Probability, Programming, Model"
$wsRes.Range("C49").Value = "Reference, Training"
$wsRes.Range("D49").Value = "Python, Statistics"
$wsRes.Range("F49").Value = "Probability for Machine Learning "
$wsRes.Range("G49").Value = "Jason Brownlee"
$wsRes.Range("H49").Value = "https://machinelearningmastery.com/probability-for-machine-learning/"

# Row 50 - Learn Git with Bitbucket Cloud
$wsRes.Range("A50").Value = "Web"
$wsRes.Range("B50").Value = "Programming"
$wsRes.Range("C50").Value = "Training"
$wsRes.Range("D50").Value = "GIT"
$wsRes.Range("F50").Value = "Learn Git with Bitbucket Cloud"
$wsRes.Range("G50").Value = "Atlassian"
$wsRes.Range("H50").Value = "https://www.atlassian.com/git"

# Hyperlinks for column H on the new rows
$wsRes.Hyperlinks.Add($wsRes.Range("H48"), "https://christophm.github.io/interpretable-ml-book/") | Out-Null
$wsRes.Hyperlinks.Add($wsRes.Range("H49"), "https://machinelearningmastery.com/probability-for-machine-learning/") | Out-Null
$wsRes.Hyperlinks.Add($wsRes.Range("H50"), "https://www.atlassian.com/git") | Out-Null

# Restore the hyperlink cell style (style index 2 == built-in "Hyperlink" style)
$wsRes.Range("H48").Style = "Hyperlink"
$wsRes.Range("H49").Style = "Hyperlink"
$wsRes.Range("H50").Style = "Hyperlink"

# ---------------------------------------------------------------
# Sheet "Authors" (sheet2) - add 2 new author rows (38-39)
# ---------------------------------------------------------------
$wsAuth = $wb.Worksheets.Item("Authors")

# Row 38 - Christoph Molnar
$wsAuth.Range("A38").Value = "Christoph Molnar"
$wsAuth.Range("B38").Value = "https://christophm.github.io"
$wsAuth.Hyperlinks.Add($wsAuth.Range("B38"), "https://christophm.github.io") | Out-Null
$wsAuth.Range("B38").Style = "Hyperlink"

# Row 39 - Atlassian
$wsAuth.Range("A39").Value = "Atlassian"
$wsAuth.Range("B39").Value = "https://www.atlassian.com/"
$wsAuth.Hyperlinks.Add($wsAuth.Range("B39"), "https://www.atlassian.com/") | Out-Null
$wsAuth.Range("B39").Style = "Hyperlink"

# ---------------------------------------------------------------
# Refresh the view / selection so each sheet shows the cell the
# author would land on after typing the new entries (one row
# below the last new row, matching the Excel auto-advance
# behaviour seen in the diff).
# ---------------------------------------------------------------
$wsAuth.Activate()
$wsAuth.Range("A40").Select()

$wsRes.Activate()
$wsRes.Range("F50").Select()
